$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5: fill in Actual End-date (F5) and change Status (G5) from
# "In-progress" style to "Completed" style (reuse G4's format, which already
# shows the green "Completed" look used elsewhere in the sheet).
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 5).Copy() | Out-Null
$ws.Cells.Item(5, 6).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (match the other date cells)
$ws.Cells.Item(5, 6).Value = 43720      # F5 Actual End-date = 2019-09-12

$ws.Cells.Item(4, 7).Copy() | Out-Null
$ws.Cells.Item(5, 7).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(5, 7).Value = "Completed"

# ---------------------------------------------------------------------------
# Row 6: Function text changes to "Unit test"; Status becomes "Not started".
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 2).Value = "Unit test"
$ws.Cells.Item(6, 7).Value = "Not started"

# ---------------------------------------------------------------------------
# Row 7 is unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 8: Item # cleared, Function becomes the new client-pages task.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).ClearContents() | Out-Null
$ws.Cells.Item(8, 2).Value = "Add more client pages to use APIs (Add new User/Book)"

# ---------------------------------------------------------------------------
# Rows 9-12 shift up one conceptual "item" (their Item # and Function text
# move to the row above relative to the old layout) and each gains concrete
# Start-date / End-Date values. Row 13 is a brand-new row holding what used
# to be row 12's content.
# ---------------------------------------------------------------------------

# Before changing row 12's formatting, snapshot its current (all-plain,
# style s=1) look into the brand-new row 13, which should end up looking
# exactly like row 12 used to (plain borders, no date format).
$ws.Range("A12:G12").Copy() | Out-Null
$ws.Range("A13:G13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Copy the date-formatted style (from row 3's C/D cells) onto the C/D cells
# of rows 9-12 before writing their new date values, so they keep the same
# "m/d/yyyy"-style date format (s=2) used throughout the rest of the sheet.
$ws.Range("C3:D3").Copy() | Out-Null
$ws.Range("C9:D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C3:D3").Copy() | Out-Null
$ws.Range("C10:D10").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("C3:D3").Copy() | Out-Null
$ws.Range("C11:D11").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("C3:D3").Copy() | Out-Null
$ws.Range("C12:D12").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Row 9
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = " Config CI/CD"
$ws.Cells.Item(9, 3).Value = 43738
$ws.Cells.Item(9, 4).Value = 43742

# Row 10
$ws.Cells.Item(10, 1).Value = 6
$ws.Cells.Item(10, 2).Value = "Register Azure Account"
$ws.Cells.Item(10, 3).Value = 43745
$ws.Cells.Item(10, 4).Value = 43749

# Row 11
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Add Docker Registry in Azure and config Auto CI/CD"
$ws.Cells.Item(11, 3).Value = 43752
$ws.Cells.Item(11, 4).Value = 43756

# Row 12
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Add Wiki page"
$ws.Cells.Item(12, 3).Value = 43759
$ws.Cells.Item(12, 4).Value = 43763

# ---------------------------------------------------------------------------
# New row 13 holds the old "Send to line manager" entry (format was already
# snapshotted above, before row 12 was reformatted).
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = "Send to line manager"
$ws.Cells.Item(13, 3).ClearContents() | Out-Null
$ws.Cells.Item(13, 4).ClearContents() | Out-Null
$ws.Cells.Item(13, 5).ClearContents() | Out-Null
$ws.Cells.Item(13, 6).ClearContents() | Out-Null
$ws.Cells.Item(13, 7).Value = "Not started"

$excel.CutCopyMode = 0
